# Insert a new row at position 21 (shifts existing rows 21-89 down to 22-90),
# then populate the new row 21 with the new weekly record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("21:21").Insert()

$ws.Cells.Item(21, 1).Value = 4
$ws.Cells.Item(21, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(21, 3).Value = "Los Lagos"
$ws.Cells.Item(21, 4).Value = 45250
$ws.Cells.Item(21, 5).Value = 10
$ws.Cells.Item(21, 6).Value = 300000000
$ws.Cells.Item(21, 7).Value = "Espárragos"
$ws.Cells.Item(21, 8).Value = "Sin especificar"
$ws.Cells.Item(21, 9).Value = "Primera"
$ws.Cells.Item(21, 10).Value = 200
$ws.Cells.Item(21, 11).Value = 2000
$ws.Cells.Item(21, 12).Value = 2000
$ws.Cells.Item(21, 13).Value = 2000
$ws.Cells.Item(21, 14).Value = "$/kilo"
$ws.Cells.Item(21, 15).Value = "Provincia de Linares"
$ws.Cells.Item(21, 16).Value = 2000
$ws.Cells.Item(21, 17).Value = 1
$ws.Cells.Item(21, 18).Value = "Hortaliza"
